$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-12-09 01:53:58"

for ($row = 2; $row -le 23; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
